# Adds the two new students ("Antonio ADM" and "Pedro ADM") to the
# "jose pavan" / 9A roster, mirroring rows 32-33 of the target sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New data rows -----------------------------------------------------
# Clone the formatting of the last existing data row (31) down onto the
# two new rows, then overwrite with the new students' data. Using Copy
# (rather than hand-rolled Borders/NumberFormat calls) reuses the
# workbook's existing "bordered text" / "bordered centered number" cell
# styles instead of minting brand-new ones.
$ws.Range("A31:C31").Copy($ws.Range("A32:C32")) | Out-Null
$ws.Range("A31:C31").Copy($ws.Range("A33:C33")) | Out-Null

$ws.Range("A32").Value = "Antonio ADM"
$ws.Range("B32").Value = 554384356465
$ws.Range("C32").Value = 554384356465

$ws.Range("A33").Value = "Pedro ADM"
$ws.Range("B33").Value = 5543996440402
$ws.Range("C33").Value = 5543996440402

# --- Selection / view state ---------------------------------------------
# Matches the saved selection in the authored workbook: the two brand new
# rows selected as full rows (A32:XFD33, active cell A32).
$ws.Rows("32:33").Select() | Out-Null
